# Apple-logo coordinates workbook: add a "circle Number" / "Name" column pair
# in front of the existing X/Y/Radius/Comment table, populate the new Name
# values, add a few new Comment annotations, and highlight the "leaf" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new columns at the left. This shifts the old
#    X / Y / Radius / Comment columns (A:D) to C:F and preserves all of
#    their existing values/formatting (incl. the custom width that used
#    to live on column D, which now belongs to column F).
# ---------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2. Headers for the two new columns.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "circle Number"
$ws.Range("B1").Value = "Name"

# ---------------------------------------------------------------------
# 3. circle Number (A) + Name (B) values for each of the 16 rows.
# ---------------------------------------------------------------------
$names = @(
    "top left",
    "top right",
    "apple top",
    "left leaf",
    "right leaf",
    "Bite",
    "bite assist",
    "bottom",
    "middle small",
    "left 5",
    "right 5",
    "left 3",
    "right 3",
    "left biggest",
    "right biggest",
    "middle 13"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# ---------------------------------------------------------------------
# 4. New Comment (F) annotations.
# ---------------------------------------------------------------------
$ws.Range("F5").Value = "left leaf intersect "
$ws.Range("F6").Value = "right leaf intersect"
$ws.Range("F7").Value = "bite"
$ws.Range("F8").Value = "to assist drawing the bite"
$ws.Range("F15").Value = "minus bite 6"

# ---------------------------------------------------------------------
# 5. Highlight the "left leaf" / "right leaf" rows (C5:F6) with the
#    Gray, Accent 3, Lighter 60% fill. Theme colors must be applied one
#    cell at a time -- applying ThemeColor to a multi-cell range at once
#    picks the wrong underlying fill record.
# ---------------------------------------------------------------------
foreach ($addr in @("C5", "D5", "E5", "F5", "C6", "D6", "E6", "F6")) {
    $ws.Range($addr).Interior.Color = 15592941
}

# ---------------------------------------------------------------------
# 6. Column widths: circle Number / Name get a fitted width, and the
#    Comment column (now F) is widened to fit its longer text.
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 12.4987
$ws.Columns("B").ColumnWidth = 12.4987
$ws.Columns("F").ColumnWidth = 20.7214

# ---------------------------------------------------------------------
# 7. Misc sheet state.
# ---------------------------------------------------------------------
$ws.Range("B9").Select()
